$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 71; this shifts existing rows 71-107 down to 72-108
# and carries formatting from the row below into the new row.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new data record.
$ws.Cells.Item(71, 1).Value = 8
$ws.Cells.Item(71, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(71, 3).Value = "Coquimbo"
$ws.Cells.Item(71, 4).Value = 44553
$ws.Cells.Item(71, 5).Value = 4
$ws.Cells.Item(71, 6).Value = 100112044
$ws.Cells.Item(71, 7).Value = "Perejil"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 3200
$ws.Cells.Item(71, 11).Value = 2000
$ws.Cells.Item(71, 12).Value = 2500
$ws.Cells.Item(71, 13).Value = 2250
$ws.Cells.Item(71, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(71, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(71, 16).Value = 1500
$ws.Cells.Item(71, 17).Value = 1.5
$ws.Cells.Item(71, 18).Value = "Hortaliza"
